$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 6939.8
$ws.Range("I64").Value = 2499.6667
$ws.Range("K64").Value = 2499.6667
$ws.Range("M64").Value = -2251.6667

$ws.Range("H67").Value = 6939.8
$ws.Range("I67").Value = 2499.6667
$ws.Range("K67").Value = 2499.6667
$ws.Range("M67").Value = -1641.6667

$ws.Range("H70").Value = 2943.6667
$ws.Range("J70").Value = 2943.6667
$ws.Range("L70").Value = 8831.000100000001
$ws.Range("N70").Value = -9371.000100000001

$ws.Range("H73").Value = 2943.6667
$ws.Range("J73").Value = 2943.6667
$ws.Range("L73").Value = 8831.000100000001
$ws.Range("N73").Value = -10703.0001

$ws.Range("H107").Value = 1081.7931
$ws.Range("I107").Value = 1188.1818
$ws.Range("J107").Value = 747.4286
$ws.Range("K107").Value = 1188.1818
$ws.Range("L107").Value = 747.4286
$ws.Range("M107").Value = 731.8181999999999
$ws.Range("N107").Value = -4587.4286

$ws.Range("H116").Value = 7413.294
$ws.Range("J116").Value = 8777.556
$ws.Range("L116").Value = 8777.556
$ws.Range("N116").Value = -15661.556

$ws.Range("H141").Value = 3454.2104
$ws.Range("I141").Value = 3449.4707
$ws.Range("K141").Value = 10348.4121
$ws.Range("M141").Value = -5168.4121

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 4038
$ws.Range("I132").Value = 3879.7778
$ws.Range("K132").Value = 11639.3334
$ws.Range("M132").Value = -9109.3334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 13170
$ws.Range("I54").Value = 893.6667
$ws.Range("K54").Value = 893.6667
$ws.Range("M54").Value = -409.6667

$ws.Range("H86").Value = 125002240
$ws.Range("I86").Value = 250001070
$ws.Range("K86").Value = 250001070
$ws.Range("M86").Value = -249999947

$ws.Range("H89").Value = 125002240
$ws.Range("I89").Value = 250001070
$ws.Range("K89").Value = 1250005350
$ws.Range("M89").Value = -1249999734

$ws.Range("H134").Value = 5246.467
$ws.Range("I134").Value = 3872.3
$ws.Range("K134").Value = 11616.9
$ws.Range("M134").Value = -9081.900000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H21").Value = 1000
$ws.Range("J21").Value = 1000
$ws.Range("L21").Value = 1000
$ws.Range("N21").Value = -1470

$ws.Range("H29").Value = 9980
$ws.Range("J29").Value = 0
$ws.Range("L29").Value = 0
$ws.Range("N29").ClearContents()

$ws.Range("H31").Value = 1707.3214
$ws.Range("I31").Value = 1310.1052
$ws.Range("J31").Value = 2545.889
$ws.Range("K31").Value = 1310.1052
$ws.Range("L31").Value = 2545.889
$ws.Range("M31").Value = -1015.1052
$ws.Range("N31").Value = -3135.889

$ws.Range("H34").Value = 1707.3214
$ws.Range("I34").Value = 1310.1052
$ws.Range("J34").Value = 2545.889
$ws.Range("K34").Value = 1310.1052
$ws.Range("L34").Value = 2545.889
$ws.Range("M34").Value = -1108.1052
$ws.Range("N34").Value = -2949.889

$ws.Range("H58").Value = 53849284
$ws.Range("J58").Value = 45457348
$ws.Range("L58").Value = 45457348
$ws.Range("N58").Value = -45457754

$ws.Range("H62").Value = 108038.5
$ws.Range("I62").Value = 1000000
$ws.Range("J62").Value = 8931.666999999999
$ws.Range("K62").Value = 1000000
$ws.Range("L62").Value = 8931.666999999999
$ws.Range("M62").Value = -999376
$ws.Range("N62").Value = -10179.667

$ws.Range("H65").Value = 108038.5
$ws.Range("I65").Value = 1000000
$ws.Range("J65").Value = 8931.666999999999
$ws.Range("K65").Value = 5000000
$ws.Range("L65").Value = 44658.335
$ws.Range("M65").Value = -4996880
$ws.Range("N65").Value = -50898.335

$ws.Range("H86").Value = 9715.532999999999
$ws.Range("I86").Value = 7098.75
$ws.Range("J86").Value = 10667.091
$ws.Range("K86").Value = 7098.75
$ws.Range("L86").Value = 10667.091
$ws.Range("M86").Value = -5975.75
$ws.Range("N86").Value = -12913.091

$ws.Range("H89").Value = 9715.532999999999
$ws.Range("I89").Value = 7098.75
$ws.Range("J89").Value = 10667.091
$ws.Range("K89").Value = 35493.75
$ws.Range("L89").Value = 53335.455
$ws.Range("M89").Value = -29877.75
$ws.Range("N89").Value = -64567.455

$ws.Range("H132").Value = 7507.8887
$ws.Range("I132").Value = 5321.375
$ws.Range("K132").Value = 15964.125
$ws.Range("M132").Value = -13434.125

$ws.Range("H134").Value = 2066.7666
$ws.Range("I134").Value = 1970.1923
$ws.Range("K134").Value = 5910.5769
$ws.Range("M134").Value = -3375.5769

$ws.Range("H136").Value = 53849284
$ws.Range("J136").Value = 45457348
$ws.Range("L136").Value = 136372044
$ws.Range("N136").Value = -136377144

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1897
$ws.Range("J5").Value = 3047.4285
$ws.Range("L5").Value = 9142.2855
$ws.Range("N5").Value = -9366.2855

$ws.Range("H62").Value = 8536.182000000001
$ws.Range("J62").Value = 9099.777
$ws.Range("L62").Value = 27299.331
$ws.Range("N62").Value = -28671.331

$ws.Range("H65").Value = 8536.182000000001
$ws.Range("J65").Value = 9099.777
$ws.Range("L65").Value = 81897.993
$ws.Range("N65").Value = -88761.993

$ws.Range("H68").Value = 845.375
$ws.Range("I68").Value = 700
$ws.Range("J68").Value = 1087.6666
$ws.Range("K68").Value = 2100
$ws.Range("L68").Value = 3262.9998
$ws.Range("M68").Value = -1289
$ws.Range("N68").Value = -4884.9998

$ws.Range("H71").Value = 845.375
$ws.Range("I71").Value = 700
$ws.Range("J71").Value = 1087.6666
$ws.Range("K71").Value = 6300
$ws.Range("L71").Value = 9788.999400000001
$ws.Range("M71").Value = -2244
$ws.Range("N71").Value = -17900.9994

$ws.Range("H80").Value = 9125.75
$ws.Range("I80").Value = 2800
$ws.Range("K80").Value = 8400
$ws.Range("M80").Value = -7464

$ws.Range("H83").Value = 9125.75
$ws.Range("I83").Value = 2800
$ws.Range("K83").Value = 25200
$ws.Range("M83").Value = -20520

$ws.Range("H92").Value = 356.3846
$ws.Range("I92").Value = 119.25
$ws.Range("K92").Value = 357.75
$ws.Range("M92").Value = 890.25

$ws.Range("H98").Value = 215.4
$ws.Range("J98").Value = 215.4
$ws.Range("L98").Value = 646.2
$ws.Range("N98").Value = -3642.2

$ws.Range("H107").Value = 1462.7858
$ws.Range("J107").Value = 1786.2
$ws.Range("L107").Value = 5358.6
$ws.Range("N107").Value = -9198.6

$ws.Range("H122").Value = 647.4286
$ws.Range("I122").Value = 776.3333
$ws.Range("J122").Value = 620.7586
$ws.Range("K122").Value = 6986.9997
$ws.Range("L122").Value = 5586.8274
$ws.Range("M122").Value = -4536.9997
$ws.Range("N122").Value = -10486.8274

$ws.Range("H129").Value = 5598.5835
$ws.Range("I129").Value = 712.4
$ws.Range("J129").Value = 9088.714
$ws.Range("K129").Value = 2137.2
$ws.Range("L129").Value = 27266.142
$ws.Range("M129").Value = 2862.8
$ws.Range("N129").Value = -37266.142

$ws.Range("H132").Value = 2984.1428
$ws.Range("I132").Value = 1594
$ws.Range("J132").Value = 3215.8333
$ws.Range("K132").Value = 14346
$ws.Range("L132").Value = 28942.4997
$ws.Range("M132").Value = -11816
$ws.Range("N132").Value = -34002.4997

$ws.Range("H135").Value = 1897
$ws.Range("J135").Value = 3047.4285
$ws.Range("L135").Value = 27426.8565
$ws.Range("N135").Value = -32496.8565

$ws.Range("H140").Value = 1500.0834
$ws.Range("I140").Value = 967.2222
$ws.Range("K140").Value = 2901.6666
$ws.Range("M140").Value = 2278.3334

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 22793.273
$ws.Range("I55").Value = 23508.25
$ws.Range("K55").Value = 23508.25
$ws.Range("M55").Value = -23181.25

$ws.Range("H122").Value = 78411
$ws.Range("I122").Value = 86881.52
$ws.Range("J122").Value = 36058.4
$ws.Range("K122").Value = 260644.56
$ws.Range("L122").Value = 108175.2
$ws.Range("M122").Value = -258194.56
$ws.Range("N122").Value = -113075.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1229.8
$ws.Range("I16").Value = 1349.875
$ws.Range("K16").Value = 1349.875
$ws.Range("M16").Value = -1179.875

$ws.Range("H22").Value = 1715.9524
$ws.Range("I22").Value = 676.8
$ws.Range("J22").Value = 2660.6365
$ws.Range("K22").Value = 676.8
$ws.Range("L22").Value = 2660.6365
$ws.Range("M22").Value = -381.8
$ws.Range("N22").Value = -3250.6365

$ws.Range("H27").Value = 1715.9524
$ws.Range("I27").Value = 676.8
$ws.Range("J27").Value = 2660.6365
$ws.Range("K27").Value = 676.8
$ws.Range("L27").Value = 2660.6365
$ws.Range("M27").Value = -569.8
$ws.Range("N27").Value = -2874.6365

$ws.Range("H61").Value = 76942136
$ws.Range("I61").Value = 83353590
$ws.Range("K61").Value = 83353590
$ws.Range("M61").Value = -83353388

$ws.Range("H82").Value = 1370.9166
$ws.Range("I82").Value = 1037.2858
$ws.Range("K82").Value = 1037.2858
$ws.Range("M82").Value = -676.2858000000001

$ws.Range("H85").Value = 1370.9166
$ws.Range("I85").Value = 1037.2858
$ws.Range("K85").Value = 1037.2858
$ws.Range("M85").Value = 210.7141999999999

$ws.Range("H113").Value = 76942136
$ws.Range("I113").Value = 83353590
$ws.Range("K113").Value = 83353590
$ws.Range("M113").Value = -83351420

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3171.9524
$ws.Range("I132").Value = 3206.2778
$ws.Range("J132").Value = 2966
$ws.Range("K132").Value = 9618.8334
$ws.Range("L132").Value = 8898
$ws.Range("M132").Value = -7088.8334
$ws.Range("N132").Value = -13958

$ws.Range("H135").Value = 71977.8
$ws.Range("J135").Value = 71977.8
$ws.Range("L135").Value = 71977.8
$ws.Range("N135").Value = -82117.8
